$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("HSV Log OLS")
$ws1.Range("A2").Value = 0.144
$ws1.Range("B2").Value = 4.333
$ws1.Range("C2").Value = 0.808

$ws2 = $wb.Worksheets.Item("HSV PPML")
$ws2.Range("A2").Value = 0.012
$ws2.Range("B2").Value = 0.921
$ws2.Range("C2").Value = 0.801

$ws3 = $wb.Worksheets.Item("HSVT NLLSQ")
$ws3.Range("A2").Value = -0.078
$ws3.Range("B2").Value = 0.257
$ws3.Range("C2").Value = 16395.06
$ws3.Range("D2").Value = 0.141
